# Add two new columns, I ("I0") and J ("IF"), to the right of the existing
# H ("IP") column on the single worksheet, and populate them for every data
# row (2-37) using the same numeric formatting/style conventions as the
# neighboring H column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the bold/bordered header style already applied to H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows --------------------------------------------------------------
# row -> (I value, J value)
$values = @{
    2  = @(1, 4)
    3  = @(1, 3)
    4  = @(1, 3)
    5  = @(1, 2)
    6  = @(1, 3)
    7  = @(1, 6)
    8  = @(1, 5)
    9  = @(1, 4)
    10 = @(1, 5)
    11 = @(1, 5)
    12 = @(1, 7)
    13 = @(1, 5)
    14 = @(1, 4)
    15 = @(1, 4)
    16 = @(1, 3)
    17 = @(1, 6)
    18 = @(1, 5)
    19 = @(1, 7)
    20 = @(1, 6)
    21 = @(1, 7)
    22 = @(1, 5)
    23 = @(1, 7)
    24 = @(1, 8)
    25 = @(1, 6)
    26 = @(1, 6)
    27 = @(1, 6)
    28 = @(1, 5)
    29 = @(1, 5)
    30 = @(1, 7)
    31 = @(1, 4)
    32 = @(1, 5)
    33 = @(1, 4)
    34 = @(1, 4)
    35 = @(1, 3)
    36 = @(3, 4)
    37 = @(1, 2)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}

Write-Output "Added I0/IF columns"
